$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-31 Friday" "2025-02-01 Saturday"

Replace-Text "694×3=2082" "735×5=3675"
Replace-Text "965×7=6755" "353×3=1059"
Replace-Text "119×9=1071" "276×4=1104"
Replace-Text "582×5=2910" "947×8=7576"
Replace-Text "913×2=1826" "217×4=868"

Replace-Text "259×4=1036" "918×6=5508"
Replace-Text "683×5=3415" "913×5=4565"
Replace-Text "383×3=1149" "720×4=2880"
Replace-Text "851×7=5957" "660×5=3300"
Replace-Text "809×6=4854" "367×2=734"

Replace-Text "966×2=1932" "998×2=1996"
Replace-Text "917×4=3668" "546×4=2184"
Replace-Text "379×6=2274" "561×6=3366"
Replace-Text "701×9=6309" "782×6=4692"
Replace-Text "258×6=1548" "778×6=4668"

Replace-Text "773×3=2319" "514×5=2570"
Replace-Text "362×6=2172" "898×8=7184"
Replace-Text "616×3=1848" "896×6=5376"
Replace-Text "563×7=3941" "910×5=4550"
Replace-Text "676×8=5408" "559×8=4472"

Replace-Text "906×9=8154" "762×5=3810"
Replace-Text "977×7=6839" "927×5=4635"
Replace-Text "430×2=860" "212×5=1060"
Replace-Text "609×5=3045" "106×8=848"
Replace-Text "416×2=832" "333×9=2997"
